$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.466.83'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.81%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.879.56'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.24%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7188'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '240.31'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9999'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3126'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07817'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.17'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +7.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08246'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.84%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.887.26'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.88%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.296'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.40'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.12%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.497.50'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.940'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '248.77'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +4.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007891'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.32'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9991'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.976'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +7.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +9.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '163.88'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.064'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +1.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.36'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.366'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.484'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.385'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.66%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.154'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +3.35%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05282'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.947'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.204'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7234'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.55%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.677'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01865'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.64%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.235.28'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +8.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.724'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9051'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.11%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '73.80'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +5.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.101'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.99%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9998'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '103.88'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5336'
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.765'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.43%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000120'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.70%  '
$ws.Range('B48').Value = 'SynthetixNetwork'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.909'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +12.52%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.301'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.54%  '
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4334'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.00%  '
$ws.Range('B51').Value = 'Aptos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.089'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.15%  '
